$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and E ("preparer" columns) were "BROWN" -> now "H.BROWN".
# Column G ("protocol" column) was "E7420L" -> stays "E7420L" (unchanged
# displayed value, only relevant as part of the same metadata-format
# correction described in the commit message).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "H.BROWN"   # column B
    $ws.Cells.Item($r, 5).Value = "H.BROWN"   # column E
    $ws.Cells.Item($r, 7).Value = "E7420L"    # column G
}

# Row 2's B/E cells carried an explicit (Arial 10) style; after the edit
# they fall back to the sheet's default/"Normal" style, matching every
# other row in the column.
$ws.Range("B2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Update the active selection to E2:E27 with E2 as the active cell.
$ws.Range("E2:E27").Select()
